# Activator actions changed to slice
# Applies the OOXML-level changes described by the commit diff:
#   - Activators!B3 "red: 1" -> "red: 1" + newline + "green: 2"
#   - Activators!B4 "off: 1" -> "off: 1,2"
#   - Activators!B11 (value 1) cleared
#   - Activators!B12 "red: 0" cleared
#   - Activators!B13 "green: 0" cleared
#   - Shortcuts!C9 filled with same text as B9 ("leds off 49")
#   - Shortcuts rows 10-12 added: Merge Input=1 / =2 / =3
#   - Style xf used by Activators!B11 gets wrapText = true
#   - Selection/view state updates on Responses, Activators, Shortcuts
#   - Keep Shortcuts as the active tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Activators sheet
# ---------------------------------------------------------------------
$wsActivators = $wb.Worksheets.Item("Activators")

# "red: 1" -> "red: 1\ngreen: 2"
$wsActivators.Range("B3").Value = "red: 1`ngreen: 2"
$wsActivators.Rows.Item(3).RowHeight = 20.95

# "off: 1" -> "off: 1,2"
$wsActivators.Range("B4").Value = "off: 1,2"

# Clear the now-retired single-led actions
$wsActivators.Range("B11").ClearContents() | Out-Null
$wsActivators.Range("B12").ClearContents() | Out-Null
$wsActivators.Range("B13").ClearContents() | Out-Null

# Style used only by B11 now wraps text
$wsActivators.Range("B11").Style.WrapText = $true

# ---------------------------------------------------------------------
# Shortcuts sheet
# ---------------------------------------------------------------------
$wsShortcuts = $wb.Worksheets.Item("Shortcuts")

# C9 gets the same "leds off 49" text as B9
$wsShortcuts.Range("C9").Value = $wsShortcuts.Range("B9").Value2

# New rows for the per-input merge shortcuts
$wsShortcuts.Range("A10").Value = 1
$wsShortcuts.Range("B10").Value = "Merge Input=1"

$wsShortcuts.Range("A11").Value = 2
$wsShortcuts.Range("B11").Value = "Merge Input=2"

$wsShortcuts.Range("A12").Value = 3
$wsShortcuts.Range("B12").Value = "Merge Input=3"

# ---------------------------------------------------------------------
# View / selection state (applied last-to-first so the final Activate()
# leaves "Shortcuts" as the active tab, matching the original workbook)
# ---------------------------------------------------------------------
$wsResponses = $wb.Worksheets.Item("Responses")
$wsResponses.Range("B26").Select() | Out-Null

$wsActivators.Range("B6").Select() | Out-Null

$wsShortcuts.Range("A13").Select() | Out-Null
